$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold values that look numeric (e.g. "7.44")
# but must stay as plain text, matching the source inlineStr cells. Forcing a
# text NumberFormat before assigning prevents Excel from coercing them to
# numbers, then resetting the style back to Normal keeps the cell style index
# identical to the original (no stray "s" attribute / quotePrefix left behind).
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '68.192.52'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '3.835.77'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '600.16'
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").Value = '171.50'
$ws.Range("E6").Value = '  +3.06%  '
$ws.Range("D7").Value = '3.835.28'
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("E10").Value = '  +1.70%  '
$ws.Range("E11").Value = '  +2.43%  '
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '0.0000284'
$ws.Range("E13").Value = '  +14.96%  '
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = '4.477.08'
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").Value = '3.846.46'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = '68.258.44'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '18.36'
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").Value = '7.44'
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").Value = '10.86'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = '468.09'
$ws.Range("E22").Value = '  +0.98%  '
$ws.Range("D23").Value = '0.730'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  +0.95%  '
$ws.Range("D27").Value = '12.13'
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").Value = '10.42'
$ws.Range("E28").Value = '  +4.42%  '
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = '3.988.94'
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").Value = '7.74'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").Value = '2.29'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").Value = '31.01'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  +0.87%  '
$ws.Range("D36").Value = '3.797.46'
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").Value = '3.88'
$ws.Range("E37").Value = '  +19.79%  '
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").Value = '1.02'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.139'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D43").Value = '0.318'
$ws.Range("E43").Value = '  +2.36%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("E46").Value = '  +3.03%  '
$ws.Range("D47").Value = '417.28'
$ws.Range("E47").Value = '  -1.99%  '
$ws.Range("D48").Value = '0.000293'
$ws.Range("E48").Value = '  +7.35%  '
$ws.Range("D49").Value = '46.63'
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0360'
$ws.Range("E50").Value = '  +1.70%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '141.51'
$ws.Range("E51").Value = '  -1.51%  '

$priceVolRange.Style = "Normal"

